$d = $word.ActiveDocument

# --------------------------------------------------------------------
# 1) Near the end of the document there used to be two paragraphs:
#       - a bold paragraph repeating the page title
#       - an italic paragraph with the meta description
#    The bold "title repeat" paragraph is removed entirely, and the
#    italic paragraph's text is replaced with a new image prompt
#    (its italic run formatting is preserved). Do this BEFORE inserting
#    the new meta-description paragraph below, since that insertion
#    introduces a second, identical copy of the text being searched for.
# --------------------------------------------------------------------
$boldTitlePara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    $txt = $para.Range.Text.TrimEnd([char]13)
    if ($txt -eq "Play Dolphin Treasure Slot for Free - Aristocrat's Underwater Adventure" -and $para.Style.NameLocal -eq "Normal") {
        $boldTitlePara = $para
    }
}
if ($boldTitlePara -eq $null) {
    # Fall back: the bold repeat paragraph is always the second-to-last
    # paragraph in the document prior to this edit.
    $boldTitlePara = $d.Paragraphs($d.Paragraphs.Count - 1)
}
$boldTitlePara.Range.Delete()

$newPrompt = "Create a feature image for Dolphin Treasure that will capture the attention of online slot game enthusiasts. The image should be in a cartoon style and should feature a happy Maya warrior with glasses. The warrior should be shown diving into a blue ocean with dolphins and other sea creatures around them. The image should be colorful and eye-catching, showcasing the fun and adventurous theme of the game. The warrior should be depicted as if they are having the time of their life, enjoying the underwater world of Dolphin Treasure. The image should also include the title of the game in a fun and playful font, along with the logo of the game developer, Aristocrat."

$d.Content.Find.Execute("Dive into the ocean and play Dolphin Treasure slot for free. Enjoy exciting gameplay, high-quality graphics, and potential payouts with Aristocrat's online slot game.", $true, $false, $false, $false, $false, $true, 1, $false, $newPrompt, 2)

# --------------------------------------------------------------------
# 2) Insert a new "Meta description" paragraph right after the H1 title
#    paragraph at the top of the document.
# --------------------------------------------------------------------
$titlePara = $d.Paragraphs(1)
$titlePara.Range.InsertParagraphAfter()
$metaPara = $d.Paragraphs(2)

$metaRest = ": Dive into the ocean and play Dolphin Treasure slot for free. Enjoy exciting gameplay, high-quality graphics, and potential payouts with Aristocrat's online slot game."

$metaXml = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>' + $metaRest + '</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$metaPara.Range.InsertXML($metaXml)
